$wb = $excel.ActiveWorkbook

# --- binek_arac: update the grup1 discount rate (B2: 0.15 -> 0.2) and move the selection ---
$wsBinek = $wb.Worksheets.Item("binek_arac")
$wsBinek.Range("B2").Value = 0.2
$wsBinek.Range("B3").Select()

# --- LCV: move the selection to B2 (single cell instead of A1:B2) ---
$wsLcv = $wb.Worksheets.Item("LCV")
$wsLcv.Range("B2").Select()

# --- HDV: becomes the active/selected sheet, selection moves to B2 ---
$wsHdv = $wb.Worksheets.Item("HDV")
$wsHdv.Activate()
$wsHdv.Range("B2").Select()
